$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold/bordered/centered header style)
# into the two new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF), rows 2-16
$data = @{
    2  = @(7, 7)
    3  = @(7, 7)
    4  = @(9, 9)
    5  = @(4, 5)
    6  = @(9, 9)
    7  = @(5, 6)
    8  = @(4, 4)
    9  = @(4, 5)
    10 = @(7, 7)
    11 = @(8, 9)
    12 = @(4, 4)
    13 = @(5, 5)
    14 = @(9, 9)
    15 = @(3, 3)
    16 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
